# SMPTE test pattern live wallpaper - first version
# Applies: corrected C8/C9 values, new "base" row 14 + extended rows 15-26
# (height/width/pct data), a new color-hex column (I) with its own text
# number format, and removal of the stray G36 helper formula.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1. Fix the two mis-typed widths in the first SMPTE block.
# ---------------------------------------------------------------------
$ws.Range("C8").Value = 144
$ws.Range("C9").Value = 192

# ---------------------------------------------------------------------
# 2. New second data block (rows 14-26): B/C raw numbers, D/E
#    Height%/Width% formulas (row 14 is the new baseline, like row 2).
# ---------------------------------------------------------------------
$numFmt = "#\ ???/???"

$ws.Range("B14").Value = 504
$ws.Range("C14").Value = 672

$ws.Range("B15").Value = 336
$ws.Range("C15").Value = 96
$ws.Range("D15").NumberFormat = $numFmt
$ws.Range("E15").NumberFormat = $numFmt
$ws.Range("D15").Formula = "=B15/`$B`$14"
$ws.Range("E15").Formula = "=C15/`$C`$14"

$ws.Range("B16").Value = 42
$ws.Range("C16").Value = 96
$ws.Range("D16").NumberFormat = $numFmt
$ws.Range("E16").NumberFormat = $numFmt
$ws.Range("D16").Formula = "=B16/`$B`$2"
$ws.Range("E16").Formula = "=C16/`$C`$2"

$ws.Range("B17").Value = 126
$ws.Range("C17").Value = 120
$ws.Range("D17").NumberFormat = $numFmt
$ws.Range("E17").NumberFormat = $numFmt
$ws.Range("D17").Formula = "=B17/`$B`$2"
$ws.Range("E17").Formula = "=C17/`$C`$2"

$ws.Range("B18").Value = 126
$ws.Range("C18").Value = 32
$ws.Range("D18").NumberFormat = $numFmt
$ws.Range("E18").NumberFormat = $numFmt
$ws.Range("D18").Formula = "=B18/`$B`$2"
$ws.Range("E18").Formula = "=C18/`$C`$2"

$ws.Range("B19").Value = 126
$ws.Range("C19").Value = 96
$ws.Range("D19").NumberFormat = $numFmt
$ws.Range("E19").NumberFormat = $numFmt
$ws.Range("D19").Formula = "=B19/`$B`$2"
$ws.Range("E19").Formula = "=C19/`$C`$2"

$ws.Range("D20").NumberFormat = $numFmt
$ws.Range("E20").NumberFormat = $numFmt
$ws.Range("D20").Formula = "=B20/`$B`$2"
$ws.Range("E20").Formula = "=C20/`$C`$2"

$ws.Range("D21").NumberFormat = $numFmt
$ws.Range("E21").NumberFormat = $numFmt
$ws.Range("D21").Formula = "=B21/`$B`$2"
$ws.Range("E21").Formula = "=C21/`$C`$2"

$ws.Range("D22").NumberFormat = $numFmt
$ws.Range("E22").NumberFormat = $numFmt
$ws.Range("D22").Formula = "=B22/`$B`$2"
$ws.Range("E22").Formula = "=C22/`$C`$2"

$ws.Range("D23").NumberFormat = $numFmt
$ws.Range("E23").NumberFormat = $numFmt
$ws.Range("D23").Formula = "=B23/`$B`$2"
$ws.Range("E23").Formula = "=C23/`$C`$2"

$ws.Range("D24").NumberFormat = $numFmt
$ws.Range("E24").NumberFormat = $numFmt
$ws.Range("D24").Formula = "=B24/`$B`$2"
$ws.Range("E24").Formula = "=C24/`$C`$2"

$ws.Range("D25").NumberFormat = $numFmt
$ws.Range("E25").NumberFormat = $numFmt
$ws.Range("D25").Formula = "=B25/`$B`$2"
$ws.Range("E25").Formula = "=C25/`$C`$2"

$ws.Range("D26").NumberFormat = $numFmt
$ws.Range("E26").NumberFormat = $numFmt
$ws.Range("D26").Formula = "=B26/`$B`$2"
$ws.Range("E26").Formula = "=C26/`$C`$2"

# ---------------------------------------------------------------------
# 3. New column I: hex color for each swatch row, stored as text.
# ---------------------------------------------------------------------
$ws.Columns.Item(9).NumberFormat = "@"

$ws.Range("I2").Value = "C0C0C0"
$ws.Range("I3").Value = "C0C000"
$ws.Range("I4").Value = "00C0C0"
$ws.Range("I5").Value = "00C000"
$ws.Range("I6").Value = "C000C0"
$ws.Range("I7").Value = "C00000"
$ws.Range("I8").Value = "0000C0"
$ws.Range("I9").Value = "0000C0"
$ws.Range("I10").Value = "131313"
$ws.Range("I11").Value = "C000C0"
$ws.Range("I12").Value = "131313"
$ws.Range("I13").Value = "00C0C0"
$ws.Range("I14").Value = "131313"
$ws.Range("I15").Value = "C0C0C0"
$ws.Range("I16").Value = "00214C"
$ws.Range("I17").Value = "FFFFFF"
$ws.Range("I18").Value = "32006A"
$ws.Range("I19").Value = "131313"
$ws.Range("I20").Value = "090909"
$ws.Range("I21").Value = "131313"
$ws.Range("I22").Value = "1D1D1D"
$ws.Range("I23").Value = "131313"

# Rows 24-30 just carry the (blank) text style down column I.
$ws.Range("I24").Value = ""
$ws.Range("I25").Value = ""
$ws.Range("I26").Value = ""
$ws.Range("I27").Value = ""
$ws.Range("I28").Value = ""
$ws.Range("I29").Value = ""
$ws.Range("I30").Value = ""

# ---------------------------------------------------------------------
# 4. Column E now matches D's "best fit" width.
# ---------------------------------------------------------------------
$ws.Columns.Item(5).ColumnWidth = $ws.Columns.Item(4).ColumnWidth

# ---------------------------------------------------------------------
# 5. Drop the stray helper formula that duplicated G34.
# ---------------------------------------------------------------------
$ws.Range("G36").ClearContents()

# ---------------------------------------------------------------------
# 6. Update the view: scroll back to the top and select the last
#    color cell that was filled in.
# ---------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("I23").Select()
